$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format first for cells whose new values look numeric,
# so Excel keeps them stored as text (matching original inlineStr
# type) instead of auto-converting them into a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.415.58"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "3.148.07"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "606.49"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "148.10"
$ws.Range("E6").Value = "  -4.30%  "
$ws.Range("D8").Value = "3.144.84"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("D11").Value = "5.58"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "0.477"
$ws.Range("E12").Value = "  -4.39%  "
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").Value = "36.45"
$ws.Range("E14").Value = "  -5.21%  "
$ws.Range("D15").Value = "3.661.97"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "64.374.04"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "0.114"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "3.148.15"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "6.96"
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("D20").Value = "482.14"
$ws.Range("E20").Value = "  -4.27%  "
$ws.Range("D21").Value = "14.54"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").Value = "0.708"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").Value = "7.73"
$ws.Range("E23").Value = "  -2.51%  "
$ws.Range("D24").Value = "13.74"
$ws.Range("E24").Value = "  -4.33%  "
$ws.Range("D25").Value = "83.26"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "2.90"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("D28").Value = "8.52"
$ws.Range("E28").Value = "  -4.21%  "
$ws.Range("E29").Value = "  -5.28%  "
$ws.Range("E30").Value = "  -30.09%  "
$ws.Range("D31").Value = "6.90"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "2.73"
$ws.Range("E32").Value = "  -4.54%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").Value = "26.44"
$ws.Range("E34").Value = "  -5.79%  "
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  -4.51%  "
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  -5.05%  "
$ws.Range("D37").Value = "54.59"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0732"
$ws.Range("E38").Value = "  -5.10%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "3.08"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").Value = "453.98"
$ws.Range("E40").Value = "  -8.87%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.124"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0400"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").Value = "8.44"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "2.877.53"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("E45").Value = "  -7.43%  "
$ws.Range("D46").Value = "2.27"
$ws.Range("E46").Value = "  -6.72%  "
$ws.Range("D47").Value = "26.51"
$ws.Range("E47").Value = "  -5.03%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").Value = "118.56"
$ws.Range("E51").Value = "  -1.97%  "
